$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-17 21:18:24"
$ws.Range("E3").Value = "2026-02-17 21:18:26"
$ws.Range("O3").Value = "-4.0 °C"
$ws.Range("E4").Value = "2026-02-17 21:18:29"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "83%"
$ws.Range("J4").Value = "1018.6 hPa"
$ws.Range("O4").Value = "9.2 °C"
$ws.Range("E5").Value = "2026-02-17 21:18:31"
$ws.Range("E6").Value = "2026-02-17 21:18:34"
$ws.Range("J6").Value = "1018.5 hPa"
$ws.Range("E7").Value = "2026-02-17 21:18:36"
$ws.Range("J7").Value = "1018.4 hPa"
$ws.Range("E8").Value = "2026-02-17 21:18:39"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "75%"
$ws.Range("E9").Value = "2026-02-17 21:18:41"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "59%"
$ws.Range("N9").Value = "7.6 °C 20:59 TU"
$ws.Range("O9").Value = "12.4 °C"
$ws.Range("E10").Value = "2026-02-17 21:18:43"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "75%"
$ws.Range("E11").Value = "2026-02-17 21:18:46"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "52%"
$ws.Range("O11").Value = "7.3 °C"
$ws.Range("E12").Value = "2026-02-17 21:18:48"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "61%"
$ws.Range("N12").Value = "9.1 °C 20:44 TU"
$ws.Range("O12").Value = "12.6 °C"
$ws.Range("E13").Value = "2026-02-17 21:18:50"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "46%"
$ws.Range("J13").Value = "1018.0 hPa"
$ws.Range("O13").Value = "6.8 °C"
$ws.Range("E14").Value = "2026-02-17 21:18:53"
$ws.Range("E15").Value = "2026-02-17 21:18:55"
$ws.Range("N15").Value = "6.7 °C 20:58 TU"
$ws.Range("O15").Value = "12.0 °C"
$ws.Range("E16").Value = "2026-02-17 21:18:57"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "66%"
$ws.Range("E17").Value = "2026-02-17 21:19:00"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "82%"
$ws.Range("E18").Value = "2026-02-17 21:19:02"
$ws.Range("J18").Value = "1018.7 hPa"
$ws.Range("E19").Value = "2026-02-17 21:19:04"
$ws.Range("O19").Value = "7.2 °C"
$ws.Range("E20").Value = "2026-02-17 21:19:07"
$ws.Range("E21").Value = "2026-02-17 21:19:09"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "39%"
$ws.Range("J21").Value = "1017.1 hPa"
$ws.Range("O21").Value = "9.7 °C"
$ws.Range("E22").Value = "2026-02-17 21:19:11"
$ws.Range("E23").Value = "2026-02-17 21:19:14"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "70%"
$ws.Range("O23").Value = "-3.7 °C"
$ws.Range("E24").Value = "2026-02-17 21:19:16"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "65%"
$ws.Range("J24").Value = "1018.7 hPa"
$ws.Range("O24").Value = "12.7 °C"
$ws.Range("E25").Value = "2026-02-17 21:19:19"
$ws.Range("E26").Value = "2026-02-17 21:19:21"
$ws.Range("E27").Value = "2026-02-17 21:19:23"
$ws.Range("E28").Value = "2026-02-17 21:19:26"
$ws.Range("E29").Value = "2026-02-17 21:19:28"
$ws.Range("E30").Value = "2026-02-17 21:19:30"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "65%"
$ws.Range("N30").Value = "8.7 °C 20:59 TU"
$ws.Range("E31").Value = "2026-02-17 21:19:33"
$ws.Range("J31").Value = "1018.5 hPa"
$ws.Range("E32").Value = "2026-02-17 21:19:35"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "70%"
$ws.Range("O32").Value = "8.4 °C"
$ws.Range("E33").Value = "2026-02-17 21:19:37"
$ws.Range("J33").Value = "1017.4 hPa"
$ws.Range("O33").Value = "6.5 °C"
$ws.Range("E34").Value = "2026-02-17 21:19:40"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "54%"
$ws.Range("O34").Value = "1.1 °C"
$ws.Range("E35").Value = "2026-02-17 21:19:42"
$ws.Range("E36").Value = "2026-02-17 21:19:45"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "60%"
$ws.Range("N36").Value = "10.2 °C 20:59 TU"
$ws.Range("O36").Value = "12.4 °C"
$ws.Range("E37").Value = "2026-02-17 21:19:47"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "73%"
$ws.Range("J37").Value = "1019.2 hPa"
$ws.Range("E38").Value = "2026-02-17 21:19:49"
$ws.Range("E39").Value = "2026-02-17 21:19:52"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "59%"
$ws.Range("E40").Value = "2026-02-17 21:19:54"
$ws.Range("J40").Value = "1018.0 hPa"
$ws.Range("O40").Value = "9.5 °C"
$ws.Range("E41").Value = "2026-02-17 21:19:57"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "51%"
$ws.Range("O41").Value = "16.4 °C"
$ws.Range("E42").Value = "2026-02-17 21:19:59"
$ws.Range("N42").Value = "9.3 °C 20:58 TU"
$ws.Range("E43").Value = "2026-02-17 21:20:01"
$ws.Range("E44").Value = "2026-02-17 21:20:04"
$ws.Range("O44").Value = "-3.0 °C"
$ws.Range("E45").Value = "2026-02-17 21:20:06"
$ws.Range("N45").Value = "2.0 °C 20:46 TU"
$ws.Range("E46").Value = "2026-02-17 21:20:08"
$ws.Range("N46").Value = "10.3 °C 20:57 TU"
$ws.Range("O46").Value = "15.3 °C"
